$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: email-address PII mapping row ---
$emailPattern = '(["''\[{<]?)([a-zA-Z0-9_\-\.]+)@((\[[0-9]{1,3}\.[0-9]{1,3}\.[0-9]{1,3}\.)|(([a-zA-Z0-9\-]+\.)+))([a-zA-Z]{2,4}|[0-9]{1,3})(\]?)(["''\]}>]?)'
$ws.Range("A2").Value = $emailPattern
$ws.Range("B2").Value = "Lorem ..."
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "PII.email_address"
$ws.Range("E2").Value = "foo_"

# --- Row 3: health-information PII mapping row ---
$ws.Range("A3").Value = "health insurance"
$ws.Range("B3").Value = "Lorem ..."
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = "PII.personal_health_information"
$ws.Range("E3").Value = "foo_"

# --- Remove the old 4th data row entirely ---
$ws.Rows.Item(4).Delete() | Out-Null

# --- Column A is now much wider (holds the long regex pattern) and no longer auto "best fit" ---
$ws.Columns.Item(1).ColumnWidth = 41.71

# --- Final selection lands on E3, matching the last-edited cell ---
$ws.Range("E3").Select() | Out-Null
